$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55 (ALC)
$ws.Range("H55").Value = 285.1111
$ws.Range("I55").Value = 260.9091
$ws.Range("J55").Value = 323.14285
$ws.Range("K55").Value = 260.9091
$ws.Range("L55").Value = 323.14285
$ws.Range("M55").Value = -46.90910000000002
$ws.Range("N55").Value = -751.14285

# Row 125 (ALC)
$ws.Range("H125").Value = 932.25
$ws.Range("I125").Value = 631
$ws.Range("J125").Value = 1032.6666
$ws.Range("K125").Value = 5679
$ws.Range("L125").Value = 9293.999400000001
$ws.Range("M125").Value = -3219
$ws.Range("N125").Value = -14213.9994

# Row 132 (ALC)
$ws.Range("H132").Value = 2037.85
$ws.Range("I132").Value = 2038.6471
$ws.Range("J132").Value = 2033.3334
$ws.Range("K132").Value = 6115.9413
$ws.Range("L132").Value = 6100.0002
$ws.Range("M132").Value = -3585.9413
$ws.Range("N132").Value = -11160.0002

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 1512.1
$ws.Range("I2").Value = 1630.0769
$ws.Range("J2").Value = 1293
$ws.Range("K2").Value = 1630.0769
$ws.Range("L2").Value = 1293
$ws.Range("M2").Value = -1517.0769
$ws.Range("N2").Value = -1519

# Row 88 (ARM)
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

# Row 91 (ARM)
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

# Row 116 (ARM)
$ws.Range("H116").Value = 1512.1
$ws.Range("I116").Value = 1630.0769
$ws.Range("J116").Value = 1293
$ws.Range("K116").Value = 1630.0769
$ws.Range("L116").Value = 1293
$ws.Range("M116").Value = 663.9231
$ws.Range("N116").Value = -5881

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 1512.1
$ws.Range("I3").Value = 1630.0769
$ws.Range("J3").Value = 1293
$ws.Range("K3").Value = 1630.0769
$ws.Range("L3").Value = 1293
$ws.Range("M3").Value = -1516.0769
$ws.Range("N3").Value = -1521

# Row 86 (BSM)
$ws.Range("H86").Value = 2262.389
$ws.Range("J86").Value = 3740
$ws.Range("L86").Value = 3740
$ws.Range("N86").Value = -5986

# Row 89 (BSM)
$ws.Range("H89").Value = 2262.389
$ws.Range("J89").Value = 3740
$ws.Range("L89").Value = 18700
$ws.Range("N89").Value = -29932

# Row 97 (BSM)
$ws.Range("H97").Value = 2800
$ws.Range("I97").Value = 2800
$ws.Range("K97").Value = 2800
$ws.Range("M97").Value = -1809

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (CRP)
$ws.Range("H16").Value = 3314.4167
$ws.Range("I16").Value = 2898.2
$ws.Range("J16").Value = 3611.7144
$ws.Range("K16").Value = 2898.2
$ws.Range("L16").Value = 3611.7144
$ws.Range("M16").Value = -2611.2
$ws.Range("N16").Value = -4185.7144

# Row 31 (CRP)
$ws.Range("H31").Value = 2350.3225
$ws.Range("I31").Value = 2033.9546
$ws.Range("J31").Value = 3123.6667
$ws.Range("K31").Value = 2033.9546
$ws.Range("L31").Value = 3123.6667
$ws.Range("M31").Value = -1738.9546
$ws.Range("N31").Value = -3713.6667

# Row 34 (CRP)
$ws.Range("H34").Value = 2350.3225
$ws.Range("I34").Value = 2033.9546
$ws.Range("J34").Value = 3123.6667
$ws.Range("K34").Value = 2033.9546
$ws.Range("L34").Value = 3123.6667
$ws.Range("M34").Value = -1831.9546
$ws.Range("N34").Value = -3527.6667

# Row 94 (CRP)
$ws.Range("H94").Value = 4538.8237
$ws.Range("I94").Value = 1951.5714
$ws.Range("J94").Value = 6349.9
$ws.Range("K94").Value = 1951.5714
$ws.Range("L94").Value = 6349.9
$ws.Range("M94").Value = -1500.5714
$ws.Range("N94").Value = -7251.9

# Row 99 (CRP)
$ws.Range("H99").Value = 1838.25
$ws.Range("I99").Value = 1762.4615
$ws.Range("J99").Value = 2166.6667
$ws.Range("K99").Value = 1762.4615
$ws.Range("L99").Value = 2166.6667
$ws.Range("M99").Value = -264.4614999999999
$ws.Range("N99").Value = -5162.6667

# Row 107 (CRP)
$ws.Range("H107").Value = 758.9677
$ws.Range("I107").Value = 679.5294
$ws.Range("J107").Value = 855.4286
$ws.Range("K107").Value = 679.5294
$ws.Range("L107").Value = 855.4286
$ws.Range("M107").Value = 1240.4706
$ws.Range("N107").Value = -4695.4286

# Row 113 (CRP)
$ws.Range("H113").Value = 3314.4167
$ws.Range("I113").Value = 2898.2
$ws.Range("J113").Value = 3611.7144
$ws.Range("K113").Value = 2898.2
$ws.Range("L113").Value = 3611.7144
$ws.Range("M113").Value = -728.1999999999998
$ws.Range("N113").Value = -7951.7144

# Row 122 (CRP)
$ws.Range("H122").Value = 1470.0588
$ws.Range("I122").Value = 1164.1666
$ws.Range("J122").Value = 1636.909
$ws.Range("K122").Value = 3492.4998
$ws.Range("L122").Value = 4910.727000000001
$ws.Range("M122").Value = -1042.4998
$ws.Range("N122").Value = -9810.727000000001

# Row 126 (CRP)
$ws.Range("H126").Value = 1838.25
$ws.Range("I126").Value = 1762.4615
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 5287.3845
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -2817.3845
$ws.Range("N126").Value = -11440.0001

# Row 134 (CRP)
$ws.Range("H134").Value = 2344.1875
$ws.Range("I134").Value = 1549.4546
$ws.Range("K134").Value = 4648.3638
$ws.Range("M134").Value = -2113.3638

$ws = $wb.Worksheets.Item("CUL")
# Row 92 (CUL)
$ws.Range("H92").Value = 382.75
$ws.Range("I92").Value = 496
$ws.Range("J92").Value = 345
$ws.Range("K92").Value = 1488
$ws.Range("L92").Value = 1035
$ws.Range("M92").Value = -240
$ws.Range("N92").Value = -3531

# Row 122 (CUL)
$ws.Range("H122").Value = 6555.3613
$ws.Range("I122").Value = 10194.637
$ws.Range("K122").Value = 91751.73300000001
$ws.Range("M122").Value = -89301.73300000001

# Row 132 (CUL)
$ws.Range("H132").Value = 1123991.5
$ws.Range("I132").Value = 1302.6666
$ws.Range("J132").Value = 1685335.9
$ws.Range("K132").Value = 11723.9994
$ws.Range("L132").Value = 15168023.1
$ws.Range("M132").Value = -9193.999400000001
$ws.Range("N132").Value = -15173083.1

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM)
$ws.Range("H70").Value = 6715.579
$ws.Range("I70").Value = 6715.579
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 6715.579
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -6445.579
$ws.Range("N70").ClearContents()

# Row 73 (GSM)
$ws.Range("H73").Value = 6715.579
$ws.Range("I73").Value = 6715.579
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 6715.579
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -5779.579
$ws.Range("N73").ClearContents()

# Row 102 (GSM)
$ws.Range("H102").Value = 1955.36
$ws.Range("I102").Value = 1984.2858
$ws.Range("J102").Value = 1803.5
$ws.Range("K102").Value = 1984.2858
$ws.Range("L102").Value = 1803.5
$ws.Range("M102").Value = -362.2858000000001
$ws.Range("N102").Value = -5047.5

# Row 107 (GSM)
$ws.Range("H107").Value = 431.7143
$ws.Range("I107").Value = 362.77777
$ws.Range("J107").Value = 555.8
$ws.Range("K107").Value = 362.77777
$ws.Range("L107").Value = 555.8
$ws.Range("M107").Value = 1557.22223
$ws.Range("N107").Value = -4395.8

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (LTW)
$ws.Range("H46").Value = 431.25
$ws.Range("I46").Value = 341.66666
$ws.Range("J46").Value = 700
$ws.Range("K46").Value = 341.66666
$ws.Range("L46").Value = 700
$ws.Range("M46").Value = -153.66666
$ws.Range("N46").Value = -1076

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (WVR)
$ws.Range("H122").Value = 1658.8462
$ws.Range("I122").Value = 1198.5385
$ws.Range("K122").Value = 3595.6155
$ws.Range("M122").Value = -1145.6155
